# Auto-generated Excel COM-interop script to append sensor log rows
# to PIR, Humidity, and Temperature sheets (SeniorConnect_MasterLog.xlsx)

$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 131-143 ---
$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A131:F143").NumberFormat = "@"
$ws.Cells.Item(131, 1).Value = "2026-01-28"
$ws.Cells.Item(131, 2).Value = "16:47:57"
$ws.Cells.Item(131, 3).Value = "16:00"
$ws.Cells.Item(131, 4).Value = "Bathroom"
$ws.Cells.Item(131, 5).Value = "No Motion"
$ws.Cells.Item(131, 6).Value = "Inactive"
$ws.Cells.Item(132, 1).Value = "2026-01-28"
$ws.Cells.Item(132, 2).Value = "16:47:58"
$ws.Cells.Item(132, 3).Value = "16:00"
$ws.Cells.Item(132, 4).Value = "Bathroom"
$ws.Cells.Item(132, 5).Value = "No Motion"
$ws.Cells.Item(132, 6).Value = "Inactive"
$ws.Cells.Item(133, 1).Value = "2026-01-28"
$ws.Cells.Item(133, 2).Value = "16:48:02"
$ws.Cells.Item(133, 3).Value = "16:00"
$ws.Cells.Item(133, 4).Value = "Bathroom"
$ws.Cells.Item(133, 5).Value = "No Motion"
$ws.Cells.Item(133, 6).Value = "Inactive"
$ws.Cells.Item(134, 1).Value = "2026-01-28"
$ws.Cells.Item(134, 2).Value = "16:48:07"
$ws.Cells.Item(134, 3).Value = "16:00"
$ws.Cells.Item(134, 4).Value = "Bathroom"
$ws.Cells.Item(134, 5).Value = "No Motion"
$ws.Cells.Item(134, 6).Value = "Inactive"
$ws.Cells.Item(135, 1).Value = "2026-01-28"
$ws.Cells.Item(135, 2).Value = "16:48:12"
$ws.Cells.Item(135, 3).Value = "16:00"
$ws.Cells.Item(135, 4).Value = "Bathroom"
$ws.Cells.Item(135, 5).Value = "No Motion"
$ws.Cells.Item(135, 6).Value = "Inactive"
$ws.Cells.Item(136, 1).Value = "2026-01-28"
$ws.Cells.Item(136, 2).Value = "16:48:17"
$ws.Cells.Item(136, 3).Value = "16:00"
$ws.Cells.Item(136, 4).Value = "Bathroom"
$ws.Cells.Item(136, 5).Value = "No Motion"
$ws.Cells.Item(136, 6).Value = "Inactive"
$ws.Cells.Item(137, 1).Value = "2026-01-28"
$ws.Cells.Item(137, 2).Value = "16:48:22"
$ws.Cells.Item(137, 3).Value = "16:00"
$ws.Cells.Item(137, 4).Value = "Bathroom"
$ws.Cells.Item(137, 5).Value = "No Motion"
$ws.Cells.Item(137, 6).Value = "Inactive"
$ws.Cells.Item(138, 1).Value = "2026-01-28"
$ws.Cells.Item(138, 2).Value = "16:48:27"
$ws.Cells.Item(138, 3).Value = "16:00"
$ws.Cells.Item(138, 4).Value = "Bathroom"
$ws.Cells.Item(138, 5).Value = "No Motion"
$ws.Cells.Item(138, 6).Value = "Inactive"
$ws.Cells.Item(139, 1).Value = "2026-01-28"
$ws.Cells.Item(139, 2).Value = "16:48:32"
$ws.Cells.Item(139, 3).Value = "16:00"
$ws.Cells.Item(139, 4).Value = "Bathroom"
$ws.Cells.Item(139, 5).Value = "No Motion"
$ws.Cells.Item(139, 6).Value = "Inactive"
$ws.Cells.Item(140, 1).Value = "2026-01-28"
$ws.Cells.Item(140, 2).Value = "16:48:37"
$ws.Cells.Item(140, 3).Value = "16:00"
$ws.Cells.Item(140, 4).Value = "Bathroom"
$ws.Cells.Item(140, 5).Value = "No Motion"
$ws.Cells.Item(140, 6).Value = "Inactive"
$ws.Cells.Item(141, 1).Value = "2026-01-28"
$ws.Cells.Item(141, 2).Value = "16:48:42"
$ws.Cells.Item(141, 3).Value = "16:00"
$ws.Cells.Item(141, 4).Value = "Bathroom"
$ws.Cells.Item(141, 5).Value = "No Motion"
$ws.Cells.Item(141, 6).Value = "Inactive"
$ws.Cells.Item(142, 1).Value = "2026-01-28"
$ws.Cells.Item(142, 2).Value = "16:48:47"
$ws.Cells.Item(142, 3).Value = "16:00"
$ws.Cells.Item(142, 4).Value = "Bathroom"
$ws.Cells.Item(142, 5).Value = "No Motion"
$ws.Cells.Item(142, 6).Value = "Inactive"
$ws.Cells.Item(143, 1).Value = "2026-01-28"
$ws.Cells.Item(143, 2).Value = "16:48:52"
$ws.Cells.Item(143, 3).Value = "16:00"
$ws.Cells.Item(143, 4).Value = "Bathroom"
$ws.Cells.Item(143, 5).Value = "No Motion"
$ws.Cells.Item(143, 6).Value = "Inactive"

# --- Humidity sheet: append rows 128-143 ---
$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A128:F143").NumberFormat = "@"
$ws.Cells.Item(128, 1).Value = "2026-01-28"
$ws.Cells.Item(128, 2).Value = "16:47:56"
$ws.Cells.Item(128, 3).Value = "16:00"
$ws.Cells.Item(128, 4).Value = "Bathroom"
$ws.Cells.Item(128, 5).Value = "86.9%"
$ws.Cells.Item(128, 6).Value = "Active"
$ws.Cells.Item(129, 1).Value = "2026-01-28"
$ws.Cells.Item(129, 2).Value = "16:47:58"
$ws.Cells.Item(129, 3).Value = "16:00"
$ws.Cells.Item(129, 4).Value = "Bathroom"
$ws.Cells.Item(129, 5).Value = "87.7%"
$ws.Cells.Item(129, 6).Value = "Active"
$ws.Cells.Item(130, 1).Value = "2026-01-28"
$ws.Cells.Item(130, 2).Value = "16:47:59"
$ws.Cells.Item(130, 3).Value = "16:00"
$ws.Cells.Item(130, 4).Value = "Bathroom"
$ws.Cells.Item(130, 5).Value = "86.8%"
$ws.Cells.Item(130, 6).Value = "Active"
$ws.Cells.Item(131, 1).Value = "2026-01-28"
$ws.Cells.Item(131, 2).Value = "16:48:01"
$ws.Cells.Item(131, 3).Value = "16:00"
$ws.Cells.Item(131, 4).Value = "Bathroom"
$ws.Cells.Item(131, 5).Value = "87.8%"
$ws.Cells.Item(131, 6).Value = "Active"
$ws.Cells.Item(132, 1).Value = "2026-01-28"
$ws.Cells.Item(132, 2).Value = "16:48:05"
$ws.Cells.Item(132, 3).Value = "16:00"
$ws.Cells.Item(132, 4).Value = "Bathroom"
$ws.Cells.Item(132, 5).Value = "87.8%"
$ws.Cells.Item(132, 6).Value = "Active"
$ws.Cells.Item(133, 1).Value = "2026-01-28"
$ws.Cells.Item(133, 2).Value = "16:48:09"
$ws.Cells.Item(133, 3).Value = "16:00"
$ws.Cells.Item(133, 4).Value = "Bathroom"
$ws.Cells.Item(133, 5).Value = "87.9%"
$ws.Cells.Item(133, 6).Value = "Active"
$ws.Cells.Item(134, 1).Value = "2026-01-28"
$ws.Cells.Item(134, 2).Value = "16:48:13"
$ws.Cells.Item(134, 3).Value = "16:00"
$ws.Cells.Item(134, 4).Value = "Bathroom"
$ws.Cells.Item(134, 5).Value = "87.8%"
$ws.Cells.Item(134, 6).Value = "Active"
$ws.Cells.Item(135, 1).Value = "2026-01-28"
$ws.Cells.Item(135, 2).Value = "16:48:17"
$ws.Cells.Item(135, 3).Value = "16:00"
$ws.Cells.Item(135, 4).Value = "Bathroom"
$ws.Cells.Item(135, 5).Value = "87.0%"
$ws.Cells.Item(135, 6).Value = "Active"
$ws.Cells.Item(136, 1).Value = "2026-01-28"
$ws.Cells.Item(136, 2).Value = "16:48:21"
$ws.Cells.Item(136, 3).Value = "16:00"
$ws.Cells.Item(136, 4).Value = "Bathroom"
$ws.Cells.Item(136, 5).Value = "87.8%"
$ws.Cells.Item(136, 6).Value = "Active"
$ws.Cells.Item(137, 1).Value = "2026-01-28"
$ws.Cells.Item(137, 2).Value = "16:48:25"
$ws.Cells.Item(137, 3).Value = "16:00"
$ws.Cells.Item(137, 4).Value = "Bathroom"
$ws.Cells.Item(137, 5).Value = "87.8%"
$ws.Cells.Item(137, 6).Value = "Active"
$ws.Cells.Item(138, 1).Value = "2026-01-28"
$ws.Cells.Item(138, 2).Value = "16:48:29"
$ws.Cells.Item(138, 3).Value = "16:00"
$ws.Cells.Item(138, 4).Value = "Bathroom"
$ws.Cells.Item(138, 5).Value = "87.8%"
$ws.Cells.Item(138, 6).Value = "Active"
$ws.Cells.Item(139, 1).Value = "2026-01-28"
$ws.Cells.Item(139, 2).Value = "16:48:33"
$ws.Cells.Item(139, 3).Value = "16:00"
$ws.Cells.Item(139, 4).Value = "Bathroom"
$ws.Cells.Item(139, 5).Value = "87.8%"
$ws.Cells.Item(139, 6).Value = "Active"
$ws.Cells.Item(140, 1).Value = "2026-01-28"
$ws.Cells.Item(140, 2).Value = "16:48:38"
$ws.Cells.Item(140, 3).Value = "16:00"
$ws.Cells.Item(140, 4).Value = "Bathroom"
$ws.Cells.Item(140, 5).Value = "86.9%"
$ws.Cells.Item(140, 6).Value = "Active"
$ws.Cells.Item(141, 1).Value = "2026-01-28"
$ws.Cells.Item(141, 2).Value = "16:48:45"
$ws.Cells.Item(141, 3).Value = "16:00"
$ws.Cells.Item(141, 4).Value = "Bathroom"
$ws.Cells.Item(141, 5).Value = "87.9%"
$ws.Cells.Item(141, 6).Value = "Active"
$ws.Cells.Item(142, 1).Value = "2026-01-28"
$ws.Cells.Item(142, 2).Value = "16:48:49"
$ws.Cells.Item(142, 3).Value = "16:00"
$ws.Cells.Item(142, 4).Value = "Bathroom"
$ws.Cells.Item(142, 5).Value = "86.9%"
$ws.Cells.Item(142, 6).Value = "Active"
$ws.Cells.Item(143, 1).Value = "2026-01-28"
$ws.Cells.Item(143, 2).Value = "16:48:53"
$ws.Cells.Item(143, 3).Value = "16:00"
$ws.Cells.Item(143, 4).Value = "Bathroom"
$ws.Cells.Item(143, 5).Value = "87.9%"
$ws.Cells.Item(143, 6).Value = "Active"

# --- Temperature sheet: append rows 128-143 ---
$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A128:F143").NumberFormat = "@"
$ws.Cells.Item(128, 1).Value = "2026-01-28"
$ws.Cells.Item(128, 2).Value = "16:47:57"
$ws.Cells.Item(128, 3).Value = "16:00"
$ws.Cells.Item(128, 4).Value = "Bathroom"
$ws.Cells.Item(128, 5).Value = "22.9C"
$ws.Cells.Item(128, 6).Value = "Active"
$ws.Cells.Item(129, 1).Value = "2026-01-28"
$ws.Cells.Item(129, 2).Value = "16:47:58"
$ws.Cells.Item(129, 3).Value = "16:00"
$ws.Cells.Item(129, 4).Value = "Bathroom"
$ws.Cells.Item(129, 5).Value = "22.8C"
$ws.Cells.Item(129, 6).Value = "Active"
$ws.Cells.Item(130, 1).Value = "2026-01-28"
$ws.Cells.Item(130, 2).Value = "16:47:59"
$ws.Cells.Item(130, 3).Value = "16:00"
$ws.Cells.Item(130, 4).Value = "Bathroom"
$ws.Cells.Item(130, 5).Value = "22.8C"
$ws.Cells.Item(130, 6).Value = "Active"
$ws.Cells.Item(131, 1).Value = "2026-01-28"
$ws.Cells.Item(131, 2).Value = "16:48:02"
$ws.Cells.Item(131, 3).Value = "16:00"
$ws.Cells.Item(131, 4).Value = "Bathroom"
$ws.Cells.Item(131, 5).Value = "22.9C"
$ws.Cells.Item(131, 6).Value = "Active"
$ws.Cells.Item(132, 1).Value = "2026-01-28"
$ws.Cells.Item(132, 2).Value = "16:48:06"
$ws.Cells.Item(132, 3).Value = "16:00"
$ws.Cells.Item(132, 4).Value = "Bathroom"
$ws.Cells.Item(132, 5).Value = "22.9C"
$ws.Cells.Item(132, 6).Value = "Active"
$ws.Cells.Item(133, 1).Value = "2026-01-28"
$ws.Cells.Item(133, 2).Value = "16:48:10"
$ws.Cells.Item(133, 3).Value = "16:00"
$ws.Cells.Item(133, 4).Value = "Bathroom"
$ws.Cells.Item(133, 5).Value = "22.9C"
$ws.Cells.Item(133, 6).Value = "Active"
$ws.Cells.Item(134, 1).Value = "2026-01-28"
$ws.Cells.Item(134, 2).Value = "16:48:14"
$ws.Cells.Item(134, 3).Value = "16:00"
$ws.Cells.Item(134, 4).Value = "Bathroom"
$ws.Cells.Item(134, 5).Value = "22.8C"
$ws.Cells.Item(134, 6).Value = "Active"
$ws.Cells.Item(135, 1).Value = "2026-01-28"
$ws.Cells.Item(135, 2).Value = "16:48:18"
$ws.Cells.Item(135, 3).Value = "16:00"
$ws.Cells.Item(135, 4).Value = "Bathroom"
$ws.Cells.Item(135, 5).Value = "22.9C"
$ws.Cells.Item(135, 6).Value = "Active"
$ws.Cells.Item(136, 1).Value = "2026-01-28"
$ws.Cells.Item(136, 2).Value = "16:48:22"
$ws.Cells.Item(136, 3).Value = "16:00"
$ws.Cells.Item(136, 4).Value = "Bathroom"
$ws.Cells.Item(136, 5).Value = "22.8C"
$ws.Cells.Item(136, 6).Value = "Active"
$ws.Cells.Item(137, 1).Value = "2026-01-28"
$ws.Cells.Item(137, 2).Value = "16:48:26"
$ws.Cells.Item(137, 3).Value = "16:00"
$ws.Cells.Item(137, 4).Value = "Bathroom"
$ws.Cells.Item(137, 5).Value = "22.8C"
$ws.Cells.Item(137, 6).Value = "Active"
$ws.Cells.Item(138, 1).Value = "2026-01-28"
$ws.Cells.Item(138, 2).Value = "16:48:30"
$ws.Cells.Item(138, 3).Value = "16:00"
$ws.Cells.Item(138, 4).Value = "Bathroom"
$ws.Cells.Item(138, 5).Value = "22.8C"
$ws.Cells.Item(138, 6).Value = "Active"
$ws.Cells.Item(139, 1).Value = "2026-01-28"
$ws.Cells.Item(139, 2).Value = "16:48:34"
$ws.Cells.Item(139, 3).Value = "16:00"
$ws.Cells.Item(139, 4).Value = "Bathroom"
$ws.Cells.Item(139, 5).Value = "22.8C"
$ws.Cells.Item(139, 6).Value = "Active"
$ws.Cells.Item(140, 1).Value = "2026-01-28"
$ws.Cells.Item(140, 2).Value = "16:48:38"
$ws.Cells.Item(140, 3).Value = "16:00"
$ws.Cells.Item(140, 4).Value = "Bathroom"
$ws.Cells.Item(140, 5).Value = "22.8C"
$ws.Cells.Item(140, 6).Value = "Active"
$ws.Cells.Item(141, 1).Value = "2026-01-28"
$ws.Cells.Item(141, 2).Value = "16:48:46"
$ws.Cells.Item(141, 3).Value = "16:00"
$ws.Cells.Item(141, 4).Value = "Bathroom"
$ws.Cells.Item(141, 5).Value = "22.9C"
$ws.Cells.Item(141, 6).Value = "Active"
$ws.Cells.Item(142, 1).Value = "2026-01-28"
$ws.Cells.Item(142, 2).Value = "16:48:50"
$ws.Cells.Item(142, 3).Value = "16:00"
$ws.Cells.Item(142, 4).Value = "Bathroom"
$ws.Cells.Item(142, 5).Value = "22.8C"
$ws.Cells.Item(142, 6).Value = "Active"
$ws.Cells.Item(143, 1).Value = "2026-01-28"
$ws.Cells.Item(143, 2).Value = "16:48:54"
$ws.Cells.Item(143, 3).Value = "16:00"
$ws.Cells.Item(143, 4).Value = "Bathroom"
$ws.Cells.Item(143, 5).Value = "22.8C"
$ws.Cells.Item(143, 6).Value = "Active"

Write-Host "Appended rows to PIR, Humidity, and Temperature sheets."
